{"js": "// The \"Autore\" cell of the use-case header table originally listed four\n// names, each in its own paragraph:\n//   Amato Adriano\n//   Afeltra Angelo\n//   Fucile Andrea\n//   Rapa Giovanni\n// The edit trims that list down to a single remaining author paragraph:\n//   Afeltra Angelo\n// i.e. the paragraphs \"Amato Adriano\", \"Fucile Andrea\" and \"Rapa Giovanni\"\n// are removed, leaving only \"Afeltra Angelo\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst namesToRemove = [\"Amato Adriano\", \"Fucile Andrea\", \"Rapa Giovanni\"];\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (namesToRemove.includes(text)) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The \"Autore\" cell of the use-case header table originally listed four\n# names, each in its own paragraph:\n#   Amato Adriano\n#   Afeltra Angelo\n#   Fucile Andrea\n#   Rapa Giovanni\n# The edit trims that list down to a single remaining author paragraph:\n#   Afeltra Angelo\n# i.e. the paragraphs \"Amato Adriano\", \"Fucile Andrea\" and \"Rapa Giovanni\"\n# are removed, leaving only \"Afeltra Angelo\".\n\n$d = $word.ActiveDocument\n\n$namesToRemove = @(\"Amato Adriano\", \"Fucile Andrea\", \"Rapa Giovanni\")\n\n# Collect the paragraphs to remove first (a paragraph's text includes its\n# trailing paragraph mark - [char]13, or [char]13 + [char]7 for the last\n# paragraph in a table cell - so trim those off before comparing).\n$toDelete = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($namesToRemove -contains $t) {\n        $toDelete += $p\n    }\n}\n\n# Delete from last to first so removing one paragraph doesn't invalidate the\n# still-pending Range references held by the others.\nfor ($i = $toDelete.Count - 1; $i -ge 0; $i--) {\n    $toDelete[$i].Range.Delete()\n}\n"}
